{"js": "// Update the 20x5 multiplication-problem table (skip the date paragraph\n// above it) to the new set of \"AA\u00d7BB=\" problems, cell by cell, in document\n// order (row-major, matching the table's visual layout). A plain text\n// find/replace is unsafe here because a handful of \"before\" values repeat\n// (e.g. \"82\u00d726=\" appears twice, mapping to two different \"after\" values)\n// and some \"after\" values equal an unrelated cell's \"before\" value, so\n// position-based Table.getCell(row, col) addressing is used instead of\n// Range.search()/replace-by-text.\nconst newValues = [\n  [\"78\u00d782=\", \"92\u00d745=\", \"32\u00d742=\", \"36\u00d779=\", \"48\u00d729=\"],\n  [\"14\u00d770=\", \"18\u00d767=\", \"38\u00d768=\", \"90\u00d730=\", \"98\u00d722=\"],\n  [\"100\u00d764=\", \"42\u00d771=\", \"53\u00d746=\", \"32\u00d783=\", \"33\u00d751=\"],\n  [\"93\u00d767=\", \"68\u00d753=\", \"47\u00d717=\", \"50\u00d764=\", \"95\u00d773=\"],\n  [\"76\u00d746=\", \"95\u00d757=\", \"99\u00d783=\", \"46\u00d7100=\", \"64\u00d715=\"],\n  [\"11\u00d749=\", \"94\u00d753=\", \"57\u00d778=\", \"99\u00d726=\", \"58\u00d779=\"],\n  [\"51\u00d761=\", \"54\u00d736=\", \"56\u00d758=\", \"28\u00d744=\", \"43\u00d715=\"],\n  [\"39\u00d734=\", \"95\u00d790=\", \"83\u00d760=\", \"64\u00d726=\", \"40\u00d714=\"],\n  [\"41\u00d745=\", \"86\u00d747=\", \"35\u00d722=\", \"78\u00d797=\", \"82\u00d778=\"],\n  [\"98\u00d734=\", \"58\u00d747=\", \"87\u00d729=\", \"15\u00d782=\", \"37\u00d731=\"],\n  [\"15\u00d711=\", \"41\u00d743=\", \"99\u00d742=\", \"42\u00d778=\", \"86\u00d718=\"],\n  [\"70\u00d742=\", \"98\u00d759=\", \"50\u00d711=\", \"47\u00d779=\", \"55\u00d767=\"],\n  [\"37\u00d757=\", \"45\u00d742=\", \"74\u00d792=\", \"15\u00d724=\", \"56\u00d756=\"],\n  [\"13\u00d744=\", \"18\u00d799=\", \"66\u00d729=\", \"36\u00d711=\", \"31\u00d770=\"],\n  [\"92\u00d717=\", \"81\u00d711=\", \"45\u00d716=\", \"89\u00d768=\", \"87\u00d743=\"],\n  [\"46\u00d746=\", \"42\u00d759=\", \"43\u00d787=\", \"51\u00d788=\", \"71\u00d723=\"],\n  [\"87\u00d783=\", \"22\u00d720=\", \"43\u00d750=\", \"80\u00d785=\", \"82\u00d763=\"],\n  [\"15\u00d720=\", \"17\u00d793=\", \"56\u00d785=\", \"36\u00d796=\", \"30\u00d780=\"],\n  [\"14\u00d718=\", \"42\u00d780=\", \"29\u00d745=\", \"19\u00d718=\", \"81\u00d764=\"],\n  [\"19\u00d757=\", \"76\u00d776=\", \"78\u00d784=\", \"88\u00d711=\", \"55\u00d797=\"]\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(`expected ${newValues.length} rows, found ${table.rowCount}`);\n}\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 20x5 multiplication-problem table (skip the date paragraph\n# above it) to the new set of \"AA\u00d7BB=\" problems, cell by cell, in document\n# order (row-major). A plain Find/Replace across the whole document is\n# unsafe here because a handful of \"before\" values repeat (e.g. \"82\u00d726=\"\n# appears twice, mapping to two different \"after\" values) and some \"after\"\n# values equal an unrelated cell's \"before\" value, so position-based\n# Table.Cell(row, column) addressing is used instead of $d.Content.Find.\n\n$newValues = @(\n    @(\"78\u00d782=\", \"92\u00d745=\", \"32\u00d742=\", \"36\u00d779=\", \"48\u00d729=\"),\n    @(\"14\u00d770=\", \"18\u00d767=\", \"38\u00d768=\", \"90\u00d730=\", \"98\u00d722=\"),\n    @(\"100\u00d764=\", \"42\u00d771=\", \"53\u00d746=\", \"32\u00d783=\", \"33\u00d751=\"),\n    @(\"93\u00d767=\", \"68\u00d753=\", \"47\u00d717=\", \"50\u00d764=\", \"95\u00d773=\"),\n    @(\"76\u00d746=\", \"95\u00d757=\", \"99\u00d783=\", \"46\u00d7100=\", \"64\u00d715=\"),\n    @(\"11\u00d749=\", \"94\u00d753=\", \"57\u00d778=\", \"99\u00d726=\", \"58\u00d779=\"),\n    @(\"51\u00d761=\", \"54\u00d736=\", \"56\u00d758=\", \"28\u00d744=\", \"43\u00d715=\"),\n    @(\"39\u00d734=\", \"95\u00d790=\", \"83\u00d760=\", \"64\u00d726=\", \"40\u00d714=\"),\n    @(\"41\u00d745=\", \"86\u00d747=\", \"35\u00d722=\", \"78\u00d797=\", \"82\u00d778=\"),\n    @(\"98\u00d734=\", \"58\u00d747=\", \"87\u00d729=\", \"15\u00d782=\", \"37\u00d731=\"),\n    @(\"15\u00d711=\", \"41\u00d743=\", \"99\u00d742=\", \"42\u00d778=\", \"86\u00d718=\"),\n    @(\"70\u00d742=\", \"98\u00d759=\", \"50\u00d711=\", \"47\u00d779=\", \"55\u00d767=\"),\n    @(\"37\u00d757=\", \"45\u00d742=\", \"74\u00d792=\", \"15\u00d724=\", \"56\u00d756=\"),\n    @(\"13\u00d744=\", \"18\u00d799=\", \"66\u00d729=\", \"36\u00d711=\", \"31\u00d770=\"),\n    @(\"92\u00d717=\", \"81\u00d711=\", \"45\u00d716=\", \"89\u00d768=\", \"87\u00d743=\"),\n    @(\"46\u00d746=\", \"42\u00d759=\", \"43\u00d787=\", \"51\u00d788=\", \"71\u00d723=\"),\n    @(\"87\u00d783=\", \"22\u00d720=\", \"43\u00d750=\", \"80\u00d785=\", \"82\u00d763=\"),\n    @(\"15\u00d720=\", \"17\u00d793=\", \"56\u00d785=\", \"36\u00d796=\", \"30\u00d780=\"),\n    @(\"14\u00d718=\", \"42\u00d780=\", \"29\u00d745=\", \"19\u00d718=\", \"81\u00d764=\"),\n    @(\"19\u00d757=\", \"76\u00d776=\", \"78\u00d784=\", \"88\u00d711=\", \"55\u00d797=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nif ($t.Rows.Count -ne $newValues.Length) {\n    throw \"expected $($newValues.Length) rows, found $($t.Rows.Count)\"\n}\n\nfor ($r = 1; $r -le $newValues.Length; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Length; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
